$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "employees" (sheet2.xml)
# ---------------------------------------------------------------------------
$emp = $wb.Worksheets.Item("employees")

# Header: ID -> employeeID
$emp.Cells.Item(2,1).Value = "employeeID"

# Move "reportsTo" header from L2 to K2, drop L2
$emp.Cells.Item(2,11).Value = "reportsTo"
$emp.Cells.Item(2,12).ClearContents()

# Employee IDs (col A, rows 3-9): add 1000 to become 1001..1007
$emp.Cells.Item(3,1).Value = 1001
$emp.Cells.Item(4,1).Value = 1002
$emp.Cells.Item(5,1).Value = 1003
$emp.Cells.Item(6,1).Value = 1004
$emp.Cells.Item(7,1).Value = 1005
$emp.Cells.Item(8,1).Value = 1006
$emp.Cells.Item(9,1).Value = 1007

# reportsTo column: collapse old K (text) + L (number) pair into a single K value
# Row 3 (the VP, Raymond) reports to nobody -> "null"
$emp.Cells.Item(3,11).Value = "null"
$emp.Cells.Item(3,12).ClearContents()

# Rows 4-9 now report to employee 1001 (Raymond)
$emp.Cells.Item(4,11).Value = 1001
$emp.Cells.Item(4,12).ClearContents()

$emp.Cells.Item(5,11).Value = 1001
$emp.Cells.Item(5,12).ClearContents()

$emp.Cells.Item(6,11).Value = 1001
$emp.Cells.Item(6,12).ClearContents()

$emp.Cells.Item(7,11).Value = 1001
$emp.Cells.Item(7,12).ClearContents()

$emp.Cells.Item(8,11).Value = 1001
$emp.Cells.Item(8,12).ClearContents()

$emp.Cells.Item(9,11).Value = 1001
$emp.Cells.Item(9,12).ClearContents()

# F8 loses its highlight formatting
$emp.Cells.Item(8,6).ClearFormats()

# Drop the stray empty, formatted row 13
$emp.Rows.Item(13).Delete()

# Column width adjustments
$emp.Columns.Item(1).ColumnWidth = 11.95
$emp.Columns.Item(10).ColumnWidth = 13.2

# Move the active-cell selection
$emp.Range("M6").Select()

# ---------------------------------------------------------------------------
# Sheet "region" (sheet3.xml)
# ---------------------------------------------------------------------------
$region = $wb.Worksheets.Item("region")

$map = @{
    3  = 1001
    4  = 1001
    5  = 1001
    6  = 1001
    7  = 1001
    8  = 1001
    9  = 1001
    10 = 1001
    11 = 1001
    12 = 1001
    13 = 1001
    14 = 1001
    15 = 1001
    16 = 1001
    17 = 1001
    18 = 1001
    19 = 1001
    20 = 1001
    21 = 1001
    22 = 1001
    23 = 1001
    24 = 1001
    25 = 1001
    26 = 1001
    27 = 1001
    28 = 1001
    29 = 1001
    30 = 1001
    31 = 1001
    32 = 1001
    33 = 1001
    34 = 1001
    35 = 1001
    36 = 1001
    37 = 1002
    38 = 1002
    39 = 1002
    40 = 1003
    41 = 1004
    42 = 1004
    43 = 1004
    44 = 1005
    45 = 1005
    46 = 1005
    47 = 1006
    48 = 1006
    49 = 1006
    50 = 1007
    51 = 1007
    52 = 1007
}

foreach ($r in $map.Keys) {
    $region.Cells.Item($r,1).Value = $map[$r]
}

$region.Range("E4").Select()
